# Add columns I (I0) and J (IF) to the sheet, mirroring the style of the
# existing header row and filling in the per-row numeric data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1): copy formatting from H1 so I1/J1 pick up the
# same style (bold, bordered, centered) already used by the other headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data cells (rows 2-63): column I ("I0") and column J ("IF") values.
$iVals = @(9,5,8,6,6,10,8,7,7,6,3,8,6,9,12,7,6,5,3,11,6,7,6,5,6,5,7,9,6,6,5,6,7,2,7,8,8,9,6,6,9,7,9,7,7,7,8,9,7,7,10,7,5,6,11,10,8,4,1,8,4,4)
$jVals = @(9,6,8,7,6,10,9,8,8,7,3,8,7,9,12,9,7,7,6,12,7,7,7,7,7,6,8,9,7,7,6,6,7,3,7,8,9,9,7,7,9,8,9,8,8,9,8,9,8,8,11,8,6,7,11,10,9,6,2,9,4,4)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}

Write-Output "Added I0/IF columns"
